# Update cached market-board / leve profit figures across the Sheets workbook
# (scheduled runner refresh of currentAveragePrice* and computed profit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 852.38464
$ws.Range("I11").Value = 852.38464
$ws.Range("K11").Value = 852.38464
$ws.Range("M11").Value = -712.38464
$ws.Range("H62").Value = 1890.4736
$ws.Range("I62").Value = 2039.3572
$ws.Range("J62").Value = 1473.6
$ws.Range("K62").Value = 2039.3572
$ws.Range("L62").Value = 1473.6
$ws.Range("M62").Value = -1415.3572
$ws.Range("N62").Value = -2721.6
$ws.Range("H65").Value = 1890.4736
$ws.Range("I65").Value = 2039.3572
$ws.Range("J65").Value = 1473.6
$ws.Range("K65").Value = 10196.786
$ws.Range("L65").Value = 7368
$ws.Range("M65").Value = -7076.786
$ws.Range("N65").Value = -13608
$ws.Range("H80").Value = 7786.75
$ws.Range("I80").Value = 930.3333
$ws.Range("J80").Value = 11900.6
$ws.Range("K80").Value = 2790.9999
$ws.Range("L80").Value = 35701.8
$ws.Range("M80").Value = -1792.9999
$ws.Range("N80").Value = -37697.8
$ws.Range("H83").Value = 7786.75
$ws.Range("I83").Value = 930.3333
$ws.Range("J83").Value = 11900.6
$ws.Range("K83").Value = 8372.9997
$ws.Range("L83").Value = 107105.4
$ws.Range("M83").Value = -3380.9997
$ws.Range("N83").Value = -117089.4
$ws.Range("H88").Value = 2189.7273
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2189.7273
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2189.7273
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -3001.7273
$ws.Range("H91").Value = 2189.7273
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2189.7273
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2189.7273
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -4997.7273
$ws.Range("H98").Value = 6368.316
$ws.Range("I98").Value = 3999.8823
$ws.Range("J98").Value = 26500
$ws.Range("K98").Value = 3999.8823
$ws.Range("L98").Value = 26500
$ws.Range("M98").Value = -2501.8823
$ws.Range("N98").Value = -29496
$ws.Range("H122").Value = 6368.316
$ws.Range("I122").Value = 3999.8823
$ws.Range("J122").Value = 26500
$ws.Range("K122").Value = 11999.6469
$ws.Range("L122").Value = 79500
$ws.Range("M122").Value = -9549.6469
$ws.Range("N122").Value = -84400
$ws.Range("H129").Value = 1108.7451
$ws.Range("I129").Value = 497.75
$ws.Range("J129").Value = 1160.7446
$ws.Range("K129").Value = 1493.25
$ws.Range("L129").Value = 3482.2338
$ws.Range("M129").Value = 3506.75
$ws.Range("N129").Value = -13482.2338
$ws.Range("H132").Value = 1817.8868
$ws.Range("I132").Value = 1584.4
$ws.Range("J132").Value = 2271.889
$ws.Range("K132").Value = 4753.200000000001
$ws.Range("L132").Value = 6815.667
$ws.Range("M132").Value = -2223.200000000001
$ws.Range("N132").Value = -11875.667
$ws.Range("H137").Value = 1276.4324
$ws.Range("I137").Value = 1230.5938
$ws.Range("J137").Value = 1569.8
$ws.Range("K137").Value = 3691.7814
$ws.Range("L137").Value = 4709.4
$ws.Range("M137").Value = -1141.7814
$ws.Range("N137").Value = -9809.4
$ws.Range("H138").Value = 2611.0168
$ws.Range("I138").Value = 1605.9048
$ws.Range("J138").Value = 3166.4736
$ws.Range("K138").Value = 4817.7144
$ws.Range("L138").Value = 9499.4208
$ws.Range("M138").Value = 322.2856000000002
$ws.Range("N138").Value = -19779.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11611.523
$ws.Range("I32").Value = 12936.255
$ws.Range("J32").Value = 4325.5
$ws.Range("K32").Value = 12936.255
$ws.Range("L32").Value = 4325.5
$ws.Range("M32").Value = -12649.255
$ws.Range("N32").Value = -4899.5
$ws.Range("H61").Value = 2179.3076
$ws.Range("I61").Value = 1987.4736
$ws.Range("K61").Value = 1987.4736
$ws.Range("M61").Value = -1775.4736
$ws.Range("H74").Value = 1372
$ws.Range("I74").Value = 930
$ws.Range("K74").Value = 930
$ws.Range("M74").Value = -56
$ws.Range("H77").Value = 1372
$ws.Range("I77").Value = 930
$ws.Range("K77").Value = 4650
$ws.Range("M77").Value = -282
$ws.Range("H136").Value = 2179.3076
$ws.Range("I136").Value = 1987.4736
$ws.Range("K136").Value = 5962.4208
$ws.Range("M136").Value = -3412.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2139.8667
$ws.Range("I134").Value = 1717.0303
$ws.Range("K134").Value = 5151.090899999999
$ws.Range("M134").Value = -2616.090899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1780.8064
$ws.Range("I31").Value = 1359.32
$ws.Range("K31").Value = 1359.32
$ws.Range("M31").Value = -1064.32
$ws.Range("H34").Value = 1780.8064
$ws.Range("I34").Value = 1359.32
$ws.Range("K34").Value = 1359.32
$ws.Range("M34").Value = -1157.32
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -36232
$ws.Range("H127").Value = 45000
$ws.Range("J127").Value = 45000
$ws.Range("L127").Value = 45000
$ws.Range("N127").Value = -54920
$ws.Range("H134").Value = 2006.7273
$ws.Range("I134").Value = 1425.0312
$ws.Range("J134").Value = 3557.9167
$ws.Range("K134").Value = 4275.0936
$ws.Range("L134").Value = 10673.7501
$ws.Range("M134").Value = -1740.0936
$ws.Range("N134").Value = -15743.7501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 24.083334
$ws.Range("I14").Value = 24.083334
$ws.Range("K14").Value = 72.25000199999999
$ws.Range("M14").Value = 100.749998
$ws.Range("H80").Value = 6858.2856
$ws.Range("I80").Value = 10000.5
$ws.Range("J80").Value = 2668.6667
$ws.Range("K80").Value = 30001.5
$ws.Range("L80").Value = 8006.000100000001
$ws.Range("M80").Value = -29065.5
$ws.Range("N80").Value = -9878.000100000001
$ws.Range("H83").Value = 6858.2856
$ws.Range("I83").Value = 10000.5
$ws.Range("J83").Value = 2668.6667
$ws.Range("K83").Value = 90004.5
$ws.Range("L83").Value = 24018.0003
$ws.Range("M83").Value = -85324.5
$ws.Range("N83").Value = -33378.0003
$ws.Range("H108").Value = 2917.9333
$ws.Range("I108").Value = 693.8
$ws.Range("J108").Value = 4030
$ws.Range("K108").Value = 2081.4
$ws.Range("L108").Value = 12090
$ws.Range("M108").Value = 798.6000000000004
$ws.Range("N108").Value = -17850

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = $null
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = $null
$ws.Range("H70").Value = 5814.048
$ws.Range("I70").Value = 5813.9287
$ws.Range("J70").Value = 5814.2856
$ws.Range("K70").Value = 5813.9287
$ws.Range("L70").Value = 5814.2856
$ws.Range("M70").Value = -5543.9287
$ws.Range("N70").Value = -6354.2856
$ws.Range("H73").Value = 5814.048
$ws.Range("I73").Value = 5813.9287
$ws.Range("J73").Value = 5814.2856
$ws.Range("K73").Value = 5813.9287
$ws.Range("L73").Value = 5814.2856
$ws.Range("M73").Value = -4877.9287
$ws.Range("N73").Value = -7686.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3772.1875
$ws.Range("I7").Value = 4357.143
$ws.Range("J7").Value = 3317.2222
$ws.Range("K7").Value = 4357.143
$ws.Range("L7").Value = 3317.2222
$ws.Range("M7").Value = -4245.143
$ws.Range("N7").Value = -3541.2222
$ws.Range("H126").Value = 3772.1875
$ws.Range("I126").Value = 4357.143
$ws.Range("J126").Value = 3317.2222
$ws.Range("K126").Value = 13071.429
$ws.Range("L126").Value = 9951.6666
$ws.Range("M126").Value = -10601.429
$ws.Range("N126").Value = -14891.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 58826830
$ws.Range("I122").Value = 111113624
$ws.Range("J122").Value = 4190.25
$ws.Range("K122").Value = 333340872
$ws.Range("L122").Value = 12570.75
$ws.Range("M122").Value = -333338422
$ws.Range("N122").Value = -17470.75
